$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Financements")

# Select the row-3 data range, then clear its contents (keeping cell styles),
# mirroring the user selecting A3:F3 and deleting the useless sample line.
$range = $ws.Range("A3:F3")
$range.Select()
$range.ClearContents()
